$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.2400528213269932, 0.3753181292658889, 0.2348887549825461, 0.4846532316848265, 0.4369199514236829, 14),
    @(0.273599112058131, 0.3783303919539526, 0.2266048251956944, 0.4760302776039508, 0.4054553575916278, 13),
    @(0.323052878118673, 0.3688046801716363, 0.2639536099994781, 0.5137641579552608, 0.4172513277134101, 12),
    @(0.3213959399964313, 0.3527736814977633, 0.2405711618554991, 0.4904805417705163, 0.3885915500499728, 11),
    @(0.3040077233811113, 0.3040077233811113, 0.2525516835738419, 0.5025452055027905, 0.4218095774931774, 10),
    @(0.3195075457514495, 0.3344878761007612, 0.2334815320987104, 0.4831992674856931, 0.3844749898822031, 9),
    @(0.3427370110204659, 0.3605818570417189, 0.2023703037307001, 0.4498558699524772, 0.3114971130477678, 8),
    @(0.2987280035122604, 0.2987280035122604, 0.1585806734357395, 0.3982218896993729, 0.2844280147574629, 7),
    @(0.2911926455816474, 0.3419929353797495, 0.3451110422066679, 0.587461524022355, 0.5589109611011353, 6),
    @(0.2508072428796626, 0.4114187230689154, 0.3395202396874553, 0.5826836531836596, 0.5880220729345085, 5)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
